# This workbook is a weekly price report. A new week's worth of data
# (2 rows: "Especial" and "Primera" quality grades) is being added at the
# top of the data block (just below the header row), pushing all of the
# previously recorded rows down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows right above the current row 3 (i.e. just after
# the header row / first data row). This shifts the existing data rows
# 3-30 down to rows 5-32, carrying their formatting (incl. the date number
# format on column D) along with them - mirroring exactly what the diff's
# dimension change (A1:T30 -> A1:T32) and per-row date shifts show.
$ws.Rows("3:4").Insert()

# Fill in the new row 3 - "Especial" grade for the new reporting date.
$ws.Cells.Item(3, 1).Value = 3
$ws.Cells.Item(3, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(3, 3).Value = "Coquimbo"
$ws.Cells.Item(3, 4).Value = 44699
$ws.Cells.Item(3, 5).Value = 5
$ws.Cells.Item(3, 6).Value = "Fruta"
$ws.Cells.Item(3, 7).Value = 100107
$ws.Cells.Item(3, 8).Value = "Otros"
$ws.Cells.Item(3, 9).Value = 100107001
$ws.Cells.Item(3, 10).Value = "Caqui"
$ws.Cells.Item(3, 11).Value = "Mankaki"
$ws.Cells.Item(3, 12).Value = "Especial"
$ws.Cells.Item(3, 13).Value = 56
$ws.Cells.Item(3, 14).Value = 12000
$ws.Cells.Item(3, 15).Value = 12000
$ws.Cells.Item(3, 16).Value = 12000
$ws.Cells.Item(3, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(3, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(3, 19).Value = 1200
$ws.Cells.Item(3, 20).Value = 10

# Fill in the new row 4 - "Primera" grade for the same new reporting date.
$ws.Cells.Item(4, 1).Value = 3
$ws.Cells.Item(4, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(4, 3).Value = "Coquimbo"
$ws.Cells.Item(4, 4).Value = 44699
$ws.Cells.Item(4, 5).Value = 5
$ws.Cells.Item(4, 6).Value = "Fruta"
$ws.Cells.Item(4, 7).Value = 100107
$ws.Cells.Item(4, 8).Value = "Otros"
$ws.Cells.Item(4, 9).Value = 100107001
$ws.Cells.Item(4, 10).Value = "Caqui"
$ws.Cells.Item(4, 11).Value = "Mankaki"
$ws.Cells.Item(4, 12).Value = "Primera"
$ws.Cells.Item(4, 13).Value = 60
$ws.Cells.Item(4, 14).Value = 10000
$ws.Cells.Item(4, 15).Value = 10000
$ws.Cells.Item(4, 16).Value = 10000
$ws.Cells.Item(4, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(4, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(4, 19).Value = 1000
$ws.Cells.Item(4, 20).Value = 10
